$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'67.732.90"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -1.83%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.264.32"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -1.21%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.03%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'579.59"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -1.61%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'184.35"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -0.70%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  +0.06%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.600"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -0.56%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  -5.04%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'6.57"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -1.90%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.407"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -4.00%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'3.825.92"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -1.36%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'  -0.28%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'27.41"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -6.04%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'67.799.61"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -1.75%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Value = "'  -3.45%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'3.271.96"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -1.61%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = "'  -2.89%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'  -1.64%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'398.18"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +2.31%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'  -2.79%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  -0.04%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'70.75"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -1.77%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  -1.98%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  -4.89%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  +0.29%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'9.50"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -2.84%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'1.01"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +0.62%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  -2.76%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'22.54"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -2.69%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  -5.71%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'6.91"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -3.94%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D34").Value = "'1.24"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -6.10%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'162.44"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -0.59%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  -5.89%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  -2.05%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'26.81"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +0.73%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  -3.68%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'4.49"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -2.63%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'6.32"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -4.25%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'2.668.15"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +1.13%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'40.66"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -2.45%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  -2.25%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'  -8.35%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'  -2.19%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'24.51"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -4.13%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.0274"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -4.16%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'6.32"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +0.01%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'  -1.89%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.967"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -3.32%  "
$ws.Range("E51").Style = "Normal"
